# إضافة حدث جديد في Card24
# Fill in the previously-blank tracking columns (B:K, P) on row 18 with
# the "nan" placeholder text used throughout this sheet, and append a new
# service-log entry as row 19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# --- Row 18: columns B..K and P were empty placeholders; fill with "nan" ---
$row18Cols = 2,3,4,5,6,7,8,9,10,11,16
foreach ($col in $row18Cols) {
    $ws.Cells.Item(18, $col).Value = "nan"
}

# --- Row 19: new service record ---
# Column A ("card") and M ("Event") hold numeric-looking text in this sheet,
# so force text formatting before assigning to avoid Excel auto-converting
# them to numbers.
$cardCell = $ws.Cells.Item(19, 1)
$cardCell.NumberFormat = "@"
$cardCell.Value = "24"

$ws.Cells.Item(19, 12).Value = "3\12\2025"

$eventCell = $ws.Cells.Item(19, 13)
$eventCell.NumberFormat = "@"
$eventCell.Value = "967"

$ws.Cells.Item(19, 14).Value = "تم سن السلندر(12 شوط) والدوفر (4 شوط)"
$ws.Cells.Item(19, 15).Value = "الخبير ارول"

Write-Output ("New dimension: " + $ws.UsedRange.Address())
